$d = $word.ActiveDocument

$d.Content.Find.Execute("92÷6=15, 2", $true, $false, $false, $false, $false, $true, 1, $false, "11÷9=1, 2", 2) | Out-Null
$d.Content.Find.Execute("41÷3=13, 2", $true, $false, $false, $false, $false, $true, 1, $false, "96÷6=16, 0", 2) | Out-Null
$d.Content.Find.Execute("53÷4=13, 1", $true, $false, $false, $false, $false, $true, 1, $false, "23÷6=3, 5", 2) | Out-Null
$d.Content.Find.Execute("84÷7=12, 0", $true, $false, $false, $false, $false, $true, 1, $false, "19÷4=4, 3", 2) | Out-Null
$d.Content.Find.Execute("38÷3=12, 2", $true, $false, $false, $false, $false, $true, 1, $false, "41÷4=10, 1", 2) | Out-Null
$d.Content.Find.Execute("43÷7=6, 1", $true, $false, $false, $false, $false, $true, 1, $false, "88÷2=44, 0", 2) | Out-Null
$d.Content.Find.Execute("94÷6=15, 4", $true, $false, $false, $false, $false, $true, 1, $false, "39÷9=4, 3", 2) | Out-Null
$d.Content.Find.Execute("71÷4=17, 3", $true, $false, $false, $false, $false, $true, 1, $false, "95÷8=11, 7", 2) | Out-Null
$d.Content.Find.Execute("84÷2=42, 0", $true, $false, $false, $false, $false, $true, 1, $false, "38÷3=12, 2", 2) | Out-Null
$d.Content.Find.Execute("65÷8=8, 1", $true, $false, $false, $false, $false, $true, 1, $false, "42÷6=7, 0", 2) | Out-Null
$d.Content.Find.Execute("73÷9=8, 1", $true, $false, $false, $false, $false, $true, 1, $false, "72÷5=14, 2", 2) | Out-Null
$d.Content.Find.Execute("68÷4=17, 0", $true, $false, $false, $false, $false, $true, 1, $false, "89÷5=17, 4", 2) | Out-Null
$d.Content.Find.Execute("67÷2=33, 1", $true, $false, $false, $false, $false, $true, 1, $false, "11÷4=2, 3", 2) | Out-Null
$d.Content.Find.Execute("35÷7=5, 0", $true, $false, $false, $false, $false, $true, 1, $false, "71÷5=14, 1", 2) | Out-Null
$d.Content.Find.Execute("80÷9=8, 8", $true, $false, $false, $false, $false, $true, 1, $false, "68÷2=34, 0", 2) | Out-Null
$d.Content.Find.Execute("67÷8=8, 3", $true, $false, $false, $false, $false, $true, 1, $false, "85÷9=9, 4", 2) | Out-Null
$d.Content.Find.Execute("88÷7=12, 4", $true, $false, $false, $false, $false, $true, 1, $false, "39÷3=13, 0", 2) | Out-Null
$d.Content.Find.Execute("75÷9=8, 3", $true, $false, $false, $false, $false, $true, 1, $false, "71÷8=8, 7", 2) | Out-Null
$d.Content.Find.Execute("61÷7=8, 5", $true, $false, $false, $false, $false, $true, 1, $false, "89÷6=14, 5", 2) | Out-Null
$d.Content.Find.Execute("96÷5=19, 1", $true, $false, $false, $false, $false, $true, 1, $false, "39÷5=7, 4", 2) | Out-Null
$d.Content.Find.Execute("33÷4=8, 1", $true, $false, $false, $false, $false, $true, 1, $false, "34÷9=3, 7", 2) | Out-Null
$d.Content.Find.Execute("19÷9=2, 1", $true, $false, $false, $false, $false, $true, 1, $false, "50÷8=6, 2", 2) | Out-Null
$d.Content.Find.Execute("74÷7=10, 4", $true, $false, $false, $false, $false, $true, 1, $false, "65÷2=32, 1", 2) | Out-Null
$d.Content.Find.Execute("72÷6=12, 0", $true, $false, $false, $false, $false, $true, 1, $false, "54÷3=18, 0", 2) | Out-Null
$d.Content.Find.Execute("53÷2=26, 1", $true, $false, $false, $false, $false, $true, 1, $false, "67÷6=11, 1", 2) | Out-Null
